# Add 2022-Q3 data: insert a new quarter sheet + a new row in the summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by copying the existing "2022-Q2" sheet
#    (keeps identical formatting: margins, styles, column widths, etc.) and
#    inserting it right before "2022-Q2", matching the target tab order.
# ---------------------------------------------------------------------------
$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($srcQ2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The copied sheet only has 3 data rows (rows 2-4); we need 5 (rows 2-6).
# Clone the formatting of the last existing data row into the two new rows.
$q3.Range("A4:H4").Copy()
$q3.Range("A5:H6").PasteSpecial(-4122)

# All of B:G are text-like fields (fund code, name, and numbers-as-text) in
# the source data, so force text format before writing them.
$q3.Range("B2:G6").NumberFormat = "@"

# Row 2 -> fund 002423
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "002423"
$q3.Range("C2").Value = "华宝标普美国品质消费股票（LOF）美元"
$q3.Range("D2").Value = "3.59"
$q3.Range("E2").Value = "94.37"
$q3.Range("F2").Value = "4.49"
$q3.Range("G2").Value = "0.1612"
$q3.Range("H2").Value = 3

# Row 3 -> fund 000043
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "000043"
$q3.Range("C3").Value = "嘉实美国成长股票（QDII）人民币"
$q3.Range("D3").Value = "12.41"
$q3.Range("E3").Value = "92.80"
$q3.Range("F3").Value = "1.27"
$q3.Range("G3").Value = "0.1576"
$q3.Range("H3").Value = 10

# Row 4 -> fund 000044
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "000044"
$q3.Range("C4").Value = "嘉实美国成长股票（QDII）美元现汇"
$q3.Range("D4").Value = "12.41"
$q3.Range("E4").Value = "92.80"
$q3.Range("F4").Value = "1.27"
$q3.Range("G4").Value = "0.1576"
$q3.Range("H4").Value = 10

# Row 5 -> fund 162415
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "162415"
$q3.Range("C5").Value = "华宝标普美国品质消费股票（LOF）人民币A"
$q3.Range("D5").Value = "2.86"
$q3.Range("E5").Value = "94.37"
$q3.Range("F5").Value = "4.49"
$q3.Range("G5").Value = "0.1284"
$q3.Range("H5").Value = 3

# Row 6 -> fund 009975
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "009975"
$q3.Range("C6").Value = "华宝标普美国品质消费股票（LOF）人民币C"
$q3.Range("D6").Value = "0.73"
$q3.Range("E6").Value = "94.37"
$q3.Range("F6").Value = "4.49"
$q3.Range("G6").Value = "0.0328"
$q3.Range("H6").Value = 3

# ---------------------------------------------------------------------------
# 2) Insert the new 2022-Q3 row at the top of the "总计" (summary) sheet,
#    pushing all the existing quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The inserted row picks up stray border formatting from the row above;
# reset it to match the plain data rows below (copy format from row 3).
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.64

# Column A is just the 0-based row position; rewrite it explicitly for every
# data row (2-9) so it stays correct after the insert shifted the old rows.
for ($i = 0; $i -le 7; $i++) {
    $total.Cells.Item($i + 2, 1).Value = $i
}

# Restore the originally-active sheet/tab (the new copied sheet steals the
# "selected" tab state by default).
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Host "2022-Q3 data added"
